$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the first column header (duplicate "Unnamed: 0" -> "Unnamed: 0.1")
$ws.Range("A1").Value = "Unnamed: 0.1"

# Populate the "Percentage" column (H) with the student's Marks value
$ws.Range("H2").Value = 45
$ws.Range("H3").Value = 90
$ws.Range("H4").Value = 85
$ws.Range("H5").Value = 67
$ws.Range("H6").Value = 88

# Highlight row 2 (failing mark, 45) in red - a brand-new fill/style
$ws.Range("A2:H2").Interior.Color = 255

# Keep rows 3-6 (passing marks) highlighted the same green they already had
$ws.Range("A3:H6").Interior.Color = 32768
